$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 14292490
$ws.Range("I76").Value = 55569610
$ws.Range("J76").Value = 4256.923
$ws.Range("K76").Value = 55569610
$ws.Range("L76").Value = 4256.923
$ws.Range("M76").Value = -55569295
$ws.Range("N76").Value = -4886.923
$ws.Range("H79").Value = 14292490
$ws.Range("I79").Value = 55569610
$ws.Range("J79").Value = 4256.923
$ws.Range("K79").Value = 55569610
$ws.Range("L79").Value = 4256.923
$ws.Range("M79").Value = -55568518
$ws.Range("N79").Value = -6440.923
$ws.Range("H80").Value = 9836.096
$ws.Range("I80").Value = 5335.2
$ws.Range("J80").Value = 13927.818
$ws.Range("K80").Value = 16005.6
$ws.Range("L80").Value = 41783.454
$ws.Range("M80").Value = -15007.6
$ws.Range("N80").Value = -43779.454
$ws.Range("H83").Value = 9836.096
$ws.Range("I83").Value = 5335.2
$ws.Range("J83").Value = 13927.818
$ws.Range("K83").Value = 48016.8
$ws.Range("L83").Value = 125350.362
$ws.Range("M83").Value = -43024.8
$ws.Range("N83").Value = -135334.362
$ws.Range("H113").Value = 5218.871
$ws.Range("I113").Value = 3837.5
$ws.Range("J113").Value = 5876.6665
$ws.Range("K113").Value = 3837.5
$ws.Range("L113").Value = 5876.6665
$ws.Range("M113").Value = -583.5
$ws.Range("N113").Value = -12384.6665

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1400
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("H32").Value = 2826.5698
$ws.Range("I32").Value = 2801.0366
$ws.Range("J32").Value = 3350
$ws.Range("K32").Value = 2801.0366
$ws.Range("L32").Value = 3350
$ws.Range("M32").Value = -2514.0366
$ws.Range("N32").Value = -3924
$ws.Range("H135").Value = 53701.668
$ws.Range("J135").Value = 53701.668
$ws.Range("L135").Value = 53701.668
$ws.Range("N135").Value = -63841.668
$ws.Range("N16").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1001.2941
$ws.Range("I22").Value = 1134
$ws.Range("J22").Value = 570
$ws.Range("K22").Value = 1134
$ws.Range("L22").Value = 570
$ws.Range("M22").Value = -784
$ws.Range("N22").Value = -1270
$ws.Range("H58").Value = 4487.525
$ws.Range("I58").Value = 7707.6
$ws.Range("J58").Value = 2555.48
$ws.Range("K58").Value = 7707.6
$ws.Range("L58").Value = 2555.48
$ws.Range("M58").Value = -7504.6
$ws.Range("N58").Value = -2961.48
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41372
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126864
$ws.Range("H102").Value = 30000
$ws.Range("J102").Value = 30000
$ws.Range("L102").Value = 30000
$ws.Range("N102").Value = -34868
$ws.Range("H105").Value = 1226.6666
$ws.Range("I105").Value = 948
$ws.Range("J105").Value = 1923.3334
$ws.Range("K105").Value = 948
$ws.Range("L105").Value = 1923.3334
$ws.Range("M105").Value = 799
$ws.Range("N105").Value = -5417.3334
$ws.Range("H122").Value = 3478.5715
$ws.Range("I122").Value = 3922.6667
$ws.Range("J122").Value = 814
$ws.Range("K122").Value = 11768.0001
$ws.Range("L122").Value = 2442
$ws.Range("M122").Value = -9318.000100000001
$ws.Range("N122").Value = -7342
$ws.Range("H136").Value = 4487.525
$ws.Range("I136").Value = 7707.6
$ws.Range("J136").Value = 2555.48
$ws.Range("K136").Value = 23122.8
$ws.Range("L136").Value = 7666.440000000001
$ws.Range("M136").Value = -20572.8
$ws.Range("N136").Value = -12766.44

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 172885.67
$ws.Range("H90").Value = 172885.67
$ws.Range("H97").Value = 911.1429000000001
$ws.Range("I97").Value = 450
$ws.Range("J97").Value = 988
$ws.Range("K97").Value = 1350
$ws.Range("L97").Value = 2964
$ws.Range("M97").Value = -854
$ws.Range("N97").Value = -3956
$ws.Range("H107").Value = 694.0213
$ws.Range("J107").Value = 1953.6923
$ws.Range("L107").Value = 5861.0769
$ws.Range("N107").Value = -9701.0769

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 633.875
$ws.Range("I2").Value = 6.4
$ws.Range("J2").Value = 1679.6666
$ws.Range("K2").Value = 6.4
$ws.Range("L2").Value = 1679.6666
$ws.Range("M2").Value = 106.6
$ws.Range("N2").Value = -1905.6666
$ws.Range("H24").Value = 431774
$ws.Range("I24").Value = 752604.5
$ws.Range("J24").Value = 4000
$ws.Range("K24").Value = 752604.5
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = -752431.5
$ws.Range("N24").Value = -4346
$ws.Range("H102").Value = 6775
$ws.Range("I102").Value = 1866.6666
$ws.Range("J102").Value = 21500
$ws.Range("K102").Value = 1866.6666
$ws.Range("L102").Value = 21500
$ws.Range("M102").Value = -244.6666
$ws.Range("N102").Value = -24744
$ws.Range("H107").Value = 5537.9473
$ws.Range("I107").Value = 9292.817999999999
$ws.Range("J107").Value = 375
$ws.Range("K107").Value = 9292.817999999999
$ws.Range("L107").Value = 375
$ws.Range("M107").Value = -7372.817999999999
$ws.Range("N107").Value = -4215
$ws.Range("H119").Value = 44666.668
$ws.Range("J119").Value = 44666.668
$ws.Range("L119").Value = 44666.668
$ws.Range("N119").Value = -54342.668
$ws.Range("H126").Value = 3122.4
$ws.Range("I126").Value = 2831.3333
$ws.Range("J126").Value = 3559
$ws.Range("K126").Value = 8493.999899999999
$ws.Range("L126").Value = 10677
$ws.Range("M126").Value = -6023.999899999999
$ws.Range("N126").Value = -15617
$ws.Range("H132").Value = 3792149.8
$ws.Range("I132").Value = 6583360
$ws.Range("J132").Value = 4078.5715
$ws.Range("K132").Value = 19750080
$ws.Range("L132").Value = 12235.7145
$ws.Range("M132").Value = -19747550
$ws.Range("N132").Value = -17295.7145

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 588.3889
$ws.Range("I22").Value = 427.92856
$ws.Range("J22").Value = 1150
$ws.Range("K22").Value = 427.92856
$ws.Range("L22").Value = 1150
$ws.Range("M22").Value = -132.92856
$ws.Range("N22").Value = -1740
$ws.Range("H23").Value = 3100.8572
$ws.Range("I23").Value = 1951
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 1951
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -1721
$ws.Range("N23").Value = -10460
$ws.Range("H27").Value = 588.3889
$ws.Range("I27").Value = 427.92856
$ws.Range("J27").Value = 1150
$ws.Range("K27").Value = 427.92856
$ws.Range("L27").Value = 1150
$ws.Range("M27").Value = -320.92856
$ws.Range("N27").Value = -1364
$ws.Range("H55").Value = 189.23334
$ws.Range("I55").Value = 176.93333
$ws.Range("J55").Value = 201.53334
$ws.Range("K55").Value = 176.93333
$ws.Range("L55").Value = 201.53334
$ws.Range("M55").Value = -3.933330000000012
$ws.Range("N55").Value = -547.53334
$ws.Range("H99").Value = 18000
$ws.Range("I99").Value = 10000
$ws.Range("K99").Value = 10000
$ws.Range("M99").Value = -7005

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3000
$ws.Range("J15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("N15").Value = -3576
